$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# phone_number (column F) values below contain leading zeros / long digit
# strings that must stay literal text (matching the source inlineStr cells),
# so mark that range as Text before writing into it.
$ws.Range("F6:F15").NumberFormat = "@"

$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = 'evan'
$ws.Cells.Item(6, 3).Value = 'Jl. kenanga'
$ws.Cells.Item(6, 4).Value = 'Daerah Khusus Ibukota Jakarta'
$ws.Cells.Item(6, 5).Value = 'Jakarta'
$ws.Cells.Item(6, 6).Value = '08463746284'
$ws.Cells.Item(6, 7).Value = '2022-11-02 06:29:26.655375'
$ws.Cells.Item(6, 8).Value = 'SG005LTR'
$ws.Cells.Item(6, 9).Value = 3
$ws.Cells.Item(6, 10).Value = 2730000

$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = 'serena'
$ws.Cells.Item(7, 3).Value = 'jalanan'
$ws.Cells.Item(7, 4).Value = 'Daerah Khusus Ibukota Jakarta'
$ws.Cells.Item(7, 5).Value = 'Kuala lumpur'
$ws.Cells.Item(7, 6).Value = '12345'
$ws.Cells.Item(7, 7).Value = '2022-11-03 09:36:50.604001'
$ws.Cells.Item(7, 8).Value = 'CANIFIL'
$ws.Cells.Item(7, 9).Value = 2
$ws.Cells.Item(7, 10).Value = 1500000

$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = 'dsjfhkfdjankjfs'
$ws.Cells.Item(8, 3).Value = 'faldhiuagiuga'
$ws.Cells.Item(8, 4).Value = 'Kepulauan Riau'
$ws.Cells.Item(8, 5).Value = 'adfjgfjhivdk'
$ws.Cells.Item(8, 6).Value = '3246732487246'
$ws.Cells.Item(8, 7).Value = '2022-11-03 09:53:56.648287'
$ws.Cells.Item(8, 8).Value = 'FP001DUS'
$ws.Cells.Item(8, 9).Value = 10
$ws.Cells.Item(8, 10).Value = 62000000

$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = 'Clara'
$ws.Cells.Item(9, 3).Value = 'Jl ember no 12'
$ws.Cells.Item(9, 4).Value = 'Bengkulu'
$ws.Cells.Item(9, 5).Value = 'Jayapura'
$ws.Cells.Item(9, 6).Value = '098347724'
$ws.Cells.Item(9, 7).Value = '2022-11-03 09:55:27.354454'
$ws.Cells.Item(9, 8).Value = 'SG005LTR'
$ws.Cells.Item(9, 9).Value = 4
$ws.Cells.Item(9, 10).Value = 3640000

$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = 'ada dnk'
$ws.Cells.Item(10, 3).Value = 'jl djhdajkbu'
$ws.Cells.Item(10, 4).Value = 'Kepulauan Bangka Belitung'
$ws.Cells.Item(10, 5).Value = 'Bengkulu'
$ws.Cells.Item(10, 6).Value = '0846375673'
$ws.Cells.Item(10, 7).Value = '2022-11-03 09:56:06.749003'
$ws.Cells.Item(10, 8).Value = 'SG001DUS'
$ws.Cells.Item(10, 9).Value = 3
$ws.Cells.Item(10, 10).Value = 9900000

$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = 'kiki'
$ws.Cells.Item(11, 3).Value = 'jl senang'
$ws.Cells.Item(11, 4).Value = 'Jawa Timur'
$ws.Cells.Item(11, 5).Value = 'Surabaya'
$ws.Cells.Item(11, 6).Value = '084757672'
$ws.Cells.Item(11, 7).Value = '2022-11-03 09:56:33.227591'
$ws.Cells.Item(11, 8).Value = 'SG001LTR'
$ws.Cells.Item(11, 9).Value = 5
$ws.Cells.Item(11, 10).Value = 1450000

$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = 'Saartika dewi'
$ws.Cells.Item(12, 3).Value = 'Jl merdeka'
$ws.Cells.Item(12, 4).Value = 'Jawa Barat'
$ws.Cells.Item(12, 5).Value = 'Bandung'
$ws.Cells.Item(12, 6).Value = '0873645273'
$ws.Cells.Item(12, 7).Value = '2022-11-03 09:57:07.941687'
$ws.Cells.Item(12, 8).Value = 'MFFF1C1'
$ws.Cells.Item(12, 9).Value = 2
$ws.Cells.Item(12, 10).Value = 5400000

$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = 'java'
$ws.Cells.Item(13, 3).Value = 'Jl. kenanga'
$ws.Cells.Item(13, 4).Value = 'Nusa Tenggara Barat'
$ws.Cells.Item(13, 5).Value = 'adfjgfjhivdk'
$ws.Cells.Item(13, 6).Value = '3246732487246'
$ws.Cells.Item(13, 7).Value = '2022-11-03 10:28:02.955195'
$ws.Cells.Item(13, 8).Value = 'FL007KG'
$ws.Cells.Item(13, 9).Value = 1
$ws.Cells.Item(13, 10).Value = 2240000

$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = 'lala'
$ws.Cells.Item(14, 3).Value = 'lalaland'
$ws.Cells.Item(14, 4).Value = 'Papua Tengah'
$ws.Cells.Item(14, 5).Value = 'NYC'
$ws.Cells.Item(14, 6).Value = '01182734'
$ws.Cells.Item(14, 7).Value = '2022-11-03 10:35:53'
$ws.Cells.Item(14, 8).Value = 'SG001DUS'
$ws.Cells.Item(14, 9).Value = 3
$ws.Cells.Item(14, 10).Value = 9900000

$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = 'cathrine'
$ws.Cells.Item(15, 3).Value = 'faldhiuagiuga'
$ws.Cells.Item(15, 4).Value = 'Bali'
$ws.Cells.Item(15, 5).Value = 'WC'
$ws.Cells.Item(15, 6).Value = '92378374'
$ws.Cells.Item(15, 7).Value = '2022-11-03 10:37'
$ws.Cells.Item(15, 8).Value = 'SG001LTR'
$ws.Cells.Item(15, 9).Value = 4
$ws.Cells.Item(15, 10).Value = 1160000
